# Fuel_forecasting.docx edit script
$d = $word.ActiveDocument

# Common run properties fragment reused throughout the document body.
$rPr = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'
$pPr = '<w:pPr>' + $rPr + '</w:pPr>'

function New-PkgXml([string]$bodyXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------------
# 1. Insert four new plan entries (3.1 - 3.4) before the existing
#    "- poszukac opisow dzialania sieci neuronowych" paragraph.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("- poszukać opisów działania sieci neuronowych", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(1)
$insPos = $rng.Start

$newItems = @(
    "3.1 Opis sieci neuronowej",
    "3.2 Opis działania rekurencyjnej sieci neuronowej",
    "3.3 sposoby uczenia rekurencyjnej sieci neuronowej",
    "3.4 Wybrana metoda z uzasadnieniem"
)
foreach ($t in $newItems) {
    $ins = $d.Range($insPos, $insPos)
    $ins.InsertBefore([char]13)
    $ins2 = $d.Range($insPos, $insPos)
    $ins2.InsertBefore($t)
    $insPos = $insPos + $t.Length + 1
}

Write-Host "Step 1 done"

# ---------------------------------------------------------------------------
# 2. "4. Modele arima" -> split "Modele arima" run into "Modele " + proofErr
#    wrapped "arima" run (spell-check marks, no visible text change).
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("4. Modele arima", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.MoveEnd(1, 1)  # include trailing paragraph mark
$body = $pPr + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve">4. </w:t></w:r>' + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve">Modele </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r>' + $rPr + '<w:t>arima</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>'
$rng.Delete()
$rng.InsertXML((New-PkgXml ('<w:p>' + $body + '</w:p>')))

Write-Host "Step 2 done"

# ---------------------------------------------------------------------------
# 3. "4.2 opis dzialania modelu arima" -> split trailing " arima" run into
#    " " + proofErr wrapped "arima" run.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("4.2 opis działania modelu arima", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.MoveEnd(1, 1)  # include trailing paragraph mark
$body = $pPr + `
    '<w:r>' + $rPr + '<w:t>4.2 opis działania modelu</w:t></w:r>' + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r>' + $rPr + '<w:t>arima</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>'
$rng.Delete()
$rng.InsertXML((New-PkgXml ('<w:p>' + $body + '</w:p>')))

Write-Host "Step 3 done"

# ---------------------------------------------------------------------------
# 4. "4.5 model sarima" / "4.6 model auto_arima": drop the stray en-GB
#    language tag and add proofErr spell-check marks around "sarima" and
#    "auto_arima".
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("4.5 model sarima", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.MoveEnd(1, 1)  # include trailing paragraph mark
$body = $pPr + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve">4.5 model </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r>' + $rPr + '<w:t>sarima</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>'
$rng.Delete()
$rng.InsertXML((New-PkgXml ('<w:p>' + $body + '</w:p>')))

$rng = $d.Content
$rng.Find.Execute("4.6 model auto_arima", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.MoveEnd(1, 1)  # include trailing paragraph mark
$body = $pPr + `
    '<w:r>' + $rPr + '<w:t>4.</w:t></w:r>' + `
    '<w:r>' + $rPr + '<w:t>6</w:t></w:r>' + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve"> model </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r>' + $rPr + '<w:t>auto_arima</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>'
$rng.Delete()
$rng.InsertXML((New-PkgXml ('<w:p>' + $body + '</w:p>')))

Write-Host "Step 4 done"

# ---------------------------------------------------------------------------
# 5. "5.3 ... wczytanego dataframe'a" -> split trailing word into a proofErr
#    wrapped run.
# ---------------------------------------------------------------------------
$rsquo = [char]0x2019
$rng = $d.Content
$rng.Find.Execute("5.3 Opisz poszczególnych kolumn oraz zamieszczenie wczytanego dataframe${rsquo}a", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.MoveEnd(1, 1)  # include trailing paragraph mark
$body = $pPr + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve">5.3 Opisz poszczególnych kolumn oraz zamieszczenie wczytanego </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r>' + $rPr + "<w:t>dataframe${rsquo}a</w:t></w:r>" + `
    '<w:proofErr w:type="spellEnd"/>'
$rng.Delete()
$rng.InsertXML((New-PkgXml ('<w:p>' + $body + '</w:p>')))

Write-Host "Step 5 done"

# ---------------------------------------------------------------------------
# 6. "7. Porownanie rezultatu ... modelu arima ... modelu auto_arima" -> wrap
#    "arima" and "auto_arima" with proofErr spell-check marks.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("7. Porównanie rezultatu uzyskanego za pomocą sieci, modelu arima o dobranych manualnie parametrach oraz modelu auto_arima", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.MoveEnd(1, 1)  # include trailing paragraph mark
$body = $pPr + `
    '<w:r>' + $rPr + '<w:t>7.</w:t></w:r>' + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve"> Porównanie rezultatu uzyskanego za pomocą sieci, modelu </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r>' + $rPr + '<w:t>arima</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve"> o dobranych manualnie parametrach oraz modelu </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r>' + $rPr + '<w:t>auto_arima</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>'
$rng.Delete()
$rng.InsertXML((New-PkgXml ('<w:p>' + $body + '</w:p>')))

Write-Host "Step 6 done"

# ---------------------------------------------------------------------------
# 7. "8. Analiza wyników": mark the "8." run with a lastRenderedPageBreak
#    (this paragraph now starts a fresh page after the earlier insertions).
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("8. Analiza wyników", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.MoveEnd(1, 1)  # include trailing paragraph mark
$body = $pPr + `
    '<w:r>' + $rPr + '<w:lastRenderedPageBreak/><w:t>8.</w:t></w:r>' + `
    '<w:r>' + $rPr + '<w:t xml:space="preserve"> Analiza wyników</w:t></w:r>'
$rng.Delete()
$rng.InsertXML((New-PkgXml ('<w:p>' + $body + '</w:p>')))

Write-Host "Step 7 done"
